# Weekly update: insert a new price record as row 340 ("Tuna" / "Primera" /
# Region de Arica y Parinacota, week of 44504) into the Melon sheet for
# Vega Modelo de Temuco. All the following rows shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 340; everything at/after 340 shifts down
# by one (old row 340 becomes 341, ..., old row 382 becomes 383).
$ws.Rows.Item(340).Insert()

# Populate the newly inserted row 340 with the new record's data.
$ws.Cells.Item(340, 1).Value = 10
$ws.Cells.Item(340, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(340, 3).Value = "La Araucanía"
$ws.Cells.Item(340, 4).Value = 44504
$ws.Cells.Item(340, 5).Value = 9
$ws.Cells.Item(340, 6).Value = 100112027
$ws.Cells.Item(340, 7).Value = "Melón"
$ws.Cells.Item(340, 8).Value = "Tuna"
$ws.Cells.Item(340, 9).Value = "Primera"
$ws.Cells.Item(340, 10).Value = 110
$ws.Cells.Item(340, 11).Value = 1500
$ws.Cells.Item(340, 12).Value = 1500
$ws.Cells.Item(340, 13).Value = 1500
$ws.Cells.Item(340, 14).Value = "$/unidad"
$ws.Cells.Item(340, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(340, 16).Value = 1500
$ws.Cells.Item(340, 17).Value = 1
$ws.Cells.Item(340, 18).Value = "Hortaliza"
